# Apply updated estimates to "Full results" and "For plotting" sheets.
# (commit: "use available information for each outcome separately, with age filter")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Full results")
$ws2 = $wb.Worksheets.Item("For plotting")

# ----- Sheet "Full results" -----

# Row 5: Outcome = wealth, Model = NULL MODEL
$ws1.Range("C5").Value = 0.876297854499272
$ws1.Range("D5").Value = 0.123770409774205
$ws1.Range("E5").Value = 1.00006826427348
$ws1.Range("J5").Value = 0.123761961253836
$ws1.Range("K5").Value = 0.124509568960829
$ws1.Range("L5").Value = 0.0973769302925883
$ws1.Range("M5").Value = 0.104746416381763
$ws1.Range("N5").Value = 0.221886499253417

# Row 6: Outcome = wealth, Model = CONDITIONAL MODEL
$ws1.Range("F6").Value = 0.868927865337484
$ws1.Range("G6").Value = 0.124518068516095

# Row 7: Outcome = wealth, Model = COMPLETE MODEL
$ws1.Range("H7").Value = 0.771544287679496
$ws1.Range("I7").Value = 0.0974932776552064
$ws1.Range("O7").Value = 0.228508377635599

# ----- Sheet "For plotting" -----

# Row 2: Sibcorr / income
$ws2.Range("D2").Value = 0.318321444760205
$ws2.Range("E2").Value = 0.408776708575272

# Row 3: IOLIB / income
$ws2.Range("D3").Value = 0.30718922769377
$ws2.Range("E3").Value = 0.366906291537621

# Row 4: IORAD / income
$ws2.Range("D4").Value = 0.369460038025135
$ws2.Range("E4").Value = 0.450086097433489

# Row 5: Sibcorr / wealth
$ws2.Range("C5").Value = 0.123761961253836
$ws2.Range("D5").Value = 0.0401812406897797
$ws2.Range("E5").Value = 0.207342681817893

# Row 6: IOLIB / wealth
$ws2.Range("C6").Value = 0.221886499253417
$ws2.Range("D6").Value = 0.170564263286972
$ws2.Range("E6").Value = 0.273208735219862

# Row 7: IORAD / wealth
$ws2.Range("C7").Value = 0.228508377635599
$ws2.Range("D7").Value = 0.169776667591906
$ws2.Range("E7").Value = 0.287240087679292

# Row 8: Sibcorr / health_pc
$ws2.Range("D8").Value = 0.0654259339770781
$ws2.Range("E8").Value = 0.347868625047153

# Row 9: IOLIB / health_pc
$ws2.Range("D9").Value = 0.0888121318463161
$ws2.Range("E9").Value = 0.365633769608239

# Row 10: IORAD / health_pc
$ws2.Range("D10").Value = 0.0985785216316375
$ws2.Range("E10").Value = 0.372821723181926

# Row 11: Sibcorr / education
$ws2.Range("D11").Value = 0.101101015285838
$ws2.Range("E11").Value = 0.19266078372726

# Row 12: IOLIB / education
$ws2.Range("D12").Value = 0.0993945930581704
$ws2.Range("E12").Value = 0.20333565446448

# Row 13: IORAD / education
$ws2.Range("D13").Value = 0.121028809292652
$ws2.Range("E13").Value = 0.218149233250591
